# Re-sort the price rows (rows 2-13) by date ascending, keeping the rest
# of each row's data (A, B, C, E, F, G, H, I, J, K, L, Q, R, T) unchanged.
# Only columns D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado) and S (Precio $/Kg) move together with the
# row, in a new, date-sorted order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (sorted) values keyed by destination row.
$newValues = @{
    2  = @{ D = 44181; M = 30; N = 20000; O = 20000; P = 20000; S = 4000 }
    3  = @{ D = 44196; M = 56; N = 15000; O = 15000; P = 15000; S = 3000 }
    4  = @{ D = 44907; M = 45; N = 25000; O = 25000; P = 25000; S = 5000 }
    5  = @{ D = 44186; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
    6  = @{ D = 44931; M = 50; N = 18000; O = 18000; P = 18000; S = 3600 }
    7  = @{ D = 44914; M = 56; N = 23000; O = 23000; P = 23000; S = 4600 }
    8  = @{ D = 44175; M = 25; N = 20000; O = 20000; P = 20000; S = 4000 }
    9  = @{ D = 44188; M = 30; N = 15000; O = 15000; P = 15000; S = 3000 }
    10 = @{ D = 44193; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
    11 = @{ D = 44902; M = 35; N = 12000; O = 12000; P = 12000; S = 2400 }
    12 = @{ D = 44189; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
    13 = @{ D = 44179; M = 45; N = 20000; O = 20000; P = 20000; S = 4000 }
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Range("D$row").Value2 = $vals.D
    $ws.Range("M$row").Value2 = $vals.M
    $ws.Range("N$row").Value2 = $vals.N
    $ws.Range("O$row").Value2 = $vals.O
    $ws.Range("P$row").Value2 = $vals.P
    $ws.Range("S$row").Value2 = $vals.S
}
